{"js": "// Add the three new character styles. Note: the style object returned\n// directly by `addStyle()` does not reliably scope `.font` writes in\n// this runtime, so we re-fetch the style via `getStyles().getByName()`\n// before touching its font \u2014 that path applies the formatting only to\n// the style's own rPr (matching target sz is expressed in half-points,\n// so 14pt -> sz=28, 10pt -> sz=20, 9.5pt -> sz=19).\ncontext.document.addStyle(\"GaNStyle\", \"Character\");\ncontext.document.addStyle(\"GaNParagraph\", \"Character\");\ncontext.document.addStyle(\"GaNLinks\", \"Character\");\nawait context.sync();\n\nconst styles = context.document.getStyles();\nconst gaNStyle = styles.getByName(\"GaNStyle\");\nconst gaNParagraph = styles.getByName(\"GaNParagraph\");\nconst gaNLinks = styles.getByName(\"GaNLinks\");\n\ngaNStyle.font.name = \"Calibri\";\ngaNStyle.font.size = 14;\n\ngaNParagraph.font.name = \"Calibri\";\ngaNParagraph.font.size = 10;\n\ngaNLinks.font.name = \"Calibri\";\ngaNLinks.font.bold = true;\ngaNLinks.font.color = \"#000080\";\ngaNLinks.font.size = 9.5;\ngaNLinks.font.underline = \"Single\";\n\nawait context.sync();\n\nconst body = context.document.body;\n\n// Apply GaNStyle to every paragraph holding the \"V roku 2022 ...\" text\n// (it appears four times in this document).\nconst observeRanges = body.search(\n  \"V roku 2022 m\u00f4\u017eete pozorova\u0165 S\u00fahvezdie Lev: 14. \u2013 23. apr\u00edla, 14. \u2013 23. m\u00e1ja\",\n  { matchCase: true, matchWholeWord: false }\n);\nobserveRanges.load(\"items\");\nawait context.sync();\nobserveRanges.items.forEach((r) => {\n  r.style = \"GaNStyle\";\n});\n\n// Apply GaNParagraph to the \"St\u00e1vate sa s\u00fa\u010das\u0165ou ...\" paragraph.\nconst campaignRanges = body.search(\n  \"St\u00e1vate sa s\u00fa\u010das\u0165ou celosvetovej kampane Globe at Night, ktorej cie\u013eom je meranie sveteln\u00e9ho zne\u010distenia. Pozorovan\u00edm  S\u00fahvezdie Lev na no\u010dnej oblohe a porovn\u00e1van\u00edm skuto\u010dnej situ\u00e1cie s na\u0161imi mapkami sa nielen\u017ee dozviete, ako osvetlenie vo Va\u0161om okol\u00ed prispieva k sveteln\u00e9mu zne\u010disteniu, ale budete m\u00f4c\u0165 porovna\u0165 \u00farove\u0148 sveteln\u00e9ho zne\u010distenia aj s in\u00fdmi lokalitami z cel\u00e9ho sveta. Va\u0161e pozorovanie tie\u017e roz\u0161\u00edri online datab\u00e1zu dokumentuj\u00facu vidite\u013enos\u0165 no\u010dnej oblohy na na\u0161ej plan\u00e9te\",\n  { matchCase: true, matchWholeWord: false }\n);\ncampaignRanges.load(\"items\");\nawait context.sync();\ncampaignRanges.items.forEach((r) => {\n  r.style = \"GaNParagraph\";\n});\n\n// Apply GaNLinks to the \"Mapky v tomto dokumente pripravil ...\" credit line.\nconst creditRanges = body.search(\n  \"Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\",\n  { matchCase: true, matchWholeWord: false }\n);\ncreditRanges.load(\"items\");\nawait context.sync();\ncreditRanges.items.forEach((r) => {\n  r.style = \"GaNLinks\";\n});\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Add the three new character styles (values chosen so the resulting\n# OOXML half-point / BGR-encoded values come out as specified):\n#   GaNStyle:     Calibri, 14pt            -> sz=28\n#   GaNParagraph: Calibri, 10pt            -> sz=20\n#   GaNLinks:     Calibri, bold, navy,\n#                 9.5pt, single underline  -> sz=19, color=000080\n\n$sGaNStyle = $d.Styles.Add(\"GaNStyle\", 2)\n$sGaNStyle.Font.Name = \"Calibri\"\n$sGaNStyle.Font.Size = 14\n\n$sGaNParagraph = $d.Styles.Add(\"GaNParagraph\", 2)\n$sGaNParagraph.Font.Name = \"Calibri\"\n$sGaNParagraph.Font.Size = 10\n\n$sGaNLinks = $d.Styles.Add(\"GaNLinks\", 2)\n$sGaNLinks.Font.Name = \"Calibri\"\n$sGaNLinks.Font.Bold = $true\n$sGaNLinks.Font.Color = 8388608\n$sGaNLinks.Font.Size = 9.5\n$sGaNLinks.Font.Underline = 1\n\n# Apply the new character styles to the runs of the matching paragraphs\n# (the paragraph mark itself is excluded so only the text run gets the\n# rPr/rStyle, matching the author's edit).\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*V roku 2022 m\u00f4\u017eete pozorova\u0165*\") {\n        $r = $p.Range\n        $r.End = $r.End - 1\n        $r.Style = \"GaNStyle\"\n    }\n    elseif ($t -like \"*St\u00e1vate sa s\u00fa\u010das\u0165ou celosvetovej kampane*\") {\n        $r = $p.Range\n        $r.End = $r.End - 1\n        $r.Style = \"GaNParagraph\"\n    }\n    elseif ($t -like \"*Mapky v tomto dokumente pripravil*\") {\n        $r = $p.Range\n        $r.End = $r.End - 1\n        $r.Style = \"GaNLinks\"\n    }\n}\n"}
